$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "Price" cells can look numeric (e.g. "398.70", "55.811.49")
# Force them to remain text so Excel does not silently coerce them to
# real numbers (which would drop the "." thousands separators and any
# trailing zeros), matching how the original values were stored.
$priceCells = @("D2", "D3", "D5", "D6", "D10", "D12", "D13", "D14", "D15", "D16", "D18", "D19", "D20", "D21", "D23", "D24", "D26", "D27", "D28", "D29", "D30", "D31", "D33", "D34", "D40", "D41", "D44", "D46", "D48", "D50", "D51")
foreach ($cellAddr in $priceCells) {
    $ws.Range($cellAddr).NumberFormat = "@"
}

$ws.Range("D2").Value = "55.811.49"
$ws.Range("E2").Value = "  +8.28%  "
$ws.Range("D3").Value = "3.217.64"
$ws.Range("E3").Value = "  +3.62%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "398.70"
$ws.Range("E5").Value = "  +3.22%  "
$ws.Range("D6").Value = "109.66"
$ws.Range("E6").Value = "  +6.14%  "
$ws.Range("E7").Value = "  +3.03%  "
$ws.Range("E9").Value = "  +6.38%  "
$ws.Range("D10").Value = "39.26"
$ws.Range("E10").Value = "  +5.82%  "
$ws.Range("E11").Value = "  +5.22%  "
$ws.Range("D12").Value = "0.141"
$ws.Range("E12").Value = "  +1.93%  "
$ws.Range("D13").Value = "3.726.05"
$ws.Range("E13").Value = "  +3.69%  "
$ws.Range("D14").Value = "19.11"
$ws.Range("E14").Value = "  +2.84%  "
$ws.Range("D15").Value = "8.06"
$ws.Range("E15").Value = "  +2.86%  "
$ws.Range("D16").Value = "3.228.24"
$ws.Range("E16").Value = "  +4.25%  "
$ws.Range("E17").Value = "  +5.86%  "
$ws.Range("D18").Value = "10.64"
$ws.Range("E18").Value = "  -3.19%  "
$ws.Range("D19").Value = "55.739.79"
$ws.Range("E19").Value = "  +8.08%  "
$ws.Range("D20").Value = "3.37"
$ws.Range("E20").Value = "  +3.06%  "
$ws.Range("D21").Value = "13.10"
$ws.Range("E21").Value = "  +5.63%  "
$ws.Range("E22").Value = "  +5.81%  "
$ws.Range("D23").Value = "303.40"
$ws.Range("E23").Value = "  +13.72%  "
$ws.Range("D24").Value = "75.33"
$ws.Range("E24").Value = "  +7.75%  "
$ws.Range("E25").Value = "  +2.02%  "
$ws.Range("D26").Value = "8.23"
$ws.Range("E26").Value = "  +1.69%  "
$ws.Range("D27").Value = "28.25"
$ws.Range("E27").Value = "  +4.54%  "
$ws.Range("D28").Value = "7.51"
$ws.Range("E28").Value = "  +3.56%  "
$ws.Range("D29").Value = "0.174"
$ws.Range("E29").Value = "  +4.72%  "
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  -0.09%  "
$ws.Range("D31").Value = "11.37"
$ws.Range("E31").Value = "  +9.80%  "
$ws.Range("E32").Value = "  +3.28%  "
$ws.Range("D33").Value = "0.0493"
$ws.Range("E33").Value = "  +2.19%  "
$ws.Range("D34").Value = "36.19"
$ws.Range("E34").Value = "  +2.62%  "
$ws.Range("E35").Value = "  +2.27%  "
$ws.Range("E36").Value = "  +2.66%  "
$ws.Range("E37").Value = "  +23.53%  "
$ws.Range("E38").Value = "  +0.09%  "
$ws.Range("E39").Value = "  +4.25%  "
$ws.Range("D40").Value = "134.69"
$ws.Range("E40").Value = "  +4.40%  "
$ws.Range("D41").Value = "4.05"
$ws.Range("E41").Value = "  +9.56%  "
$ws.Range("E42").Value = "  +2.15%  "
$ws.Range("E43").Value = "  +2.83%  "
$ws.Range("D44").Value = "17.06"
$ws.Range("E44").Value = "  +2.96%  "
$ws.Range("E45").Value = "  -3.04%  "
$ws.Range("D46").Value = "22.32"
$ws.Range("E46").Value = "  -0.49%  "
$ws.Range("E47").Value = "  +2.17%  "
$ws.Range("D48").Value = "2.146.46"
$ws.Range("E48").Value = "  +3.78%  "
$ws.Range("D50").Value = "2.11"
$ws.Range("E50").Value = "  +44.12%  "
$ws.Range("D51").Value = "0.0364"
$ws.Range("E51").Value = "  +9.25%  "
